# Apply updated coin price / volume (and re-ranked name/link) values
# from the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (quote-prefixed) so Excel does not
# auto-convert numeric-/percentage-looking strings (e.g. "29.43", "7.16%")
# into numbers; this preserves the original inline-string cell content.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

Set-TextValue $ws.Range("E2") "0.81%"
Set-TextValue $ws.Range("D3") "29.43"
Set-TextValue $ws.Range("E3") "7.16%"
Set-TextValue $ws.Range("D4") "5.188"
Set-TextValue $ws.Range("E4") "1.46%"
Set-TextValue $ws.Range("D5") "0.05733"
Set-TextValue $ws.Range("E5") "0.86%"
Set-TextValue $ws.Range("D6") "6.567"
Set-TextValue $ws.Range("E6") "0.77%"
Set-TextValue $ws.Range("D7") "0.8583"
Set-TextValue $ws.Range("E7") "4.71%"
Set-TextValue $ws.Range("D8") "0.8676"
Set-TextValue $ws.Range("E8") "1.70%"
Set-TextValue $ws.Range("B9") "One"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D9") "0.0006013"
Set-TextValue $ws.Range("E9") "0.22%"
Set-TextValue $ws.Range("B10") "WazirX"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1367"
Set-TextValue $ws.Range("E10") "2.29%"
Set-TextValue $ws.Range("B11") "MandalaExchangeToken"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07075"
Set-TextValue $ws.Range("E11") "1.89%"
Set-TextValue $ws.Range("B12") "BitrueCoin"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.03019"
Set-TextValue $ws.Range("E12") "4.92%"
Set-TextValue $ws.Range("B13") "BitMartToken"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09385"
Set-TextValue $ws.Range("E13") "-0.10%"
Set-TextValue $ws.Range("B14") "BitForexToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001526"
Set-TextValue $ws.Range("E14") "0.98%"
Set-TextValue $ws.Range("D15") "0.006008"
Set-TextValue $ws.Range("E15") "-3.32%"
Set-TextValue $ws.Range("E16") "5,225.28%"
Set-TextValue $ws.Range("D17") "3.494"
Set-TextValue $ws.Range("E17") "-0.46%"
Set-TextValue $ws.Range("D18") "3.104"
Set-TextValue $ws.Range("E18") "3.15%"
Set-TextValue $ws.Range("D19") "2.274"
Set-TextValue $ws.Range("E19") "1.96%"
Set-TextValue $ws.Range("E20") "-0.25%"
Set-TextValue $ws.Range("D21") "0.03296"
Set-TextValue $ws.Range("E21") "2.35%"
Set-TextValue $ws.Range("E22") "1.28%"
Set-TextValue $ws.Range("D23") "3.486"
Set-TextValue $ws.Range("E23") "-2.16%"
Set-TextValue $ws.Range("D24") "0.04146"
Set-TextValue $ws.Range("E24") "3.00%"
Set-TextValue $ws.Range("E25") "0.45%"
Set-TextValue $ws.Range("D26") "0.001226"
Set-TextValue $ws.Range("E26") "0.94%"
Set-TextValue $ws.Range("D27") "0.004996"
Set-TextValue $ws.Range("E27") "11.54%"
Set-TextValue $ws.Range("E28") "2.60%"
Set-TextValue $ws.Range("D40") "0.03754"
Set-TextValue $ws.Range("E40") "0.99%"
Set-TextValue $ws.Range("D41") "0.005793"
Set-TextValue $ws.Range("E41") "67.99%"
Set-TextValue $ws.Range("E42") "1.05%"
Set-TextValue $ws.Range("E43") "-11.72%"
Set-TextValue $ws.Range("D44") "0.009500"
Set-TextValue $ws.Range("E44") "-2.16%"
Set-TextValue $ws.Range("D45") "0.00005284"
Set-TextValue $ws.Range("E45") "3.23%"
Set-TextValue $ws.Range("D47") "0.05703"
Set-TextValue $ws.Range("E47") "-43.53%"
Set-TextValue $ws.Range("D48") "0.002278"
Set-TextValue $ws.Range("E48") "-9.49%"
